$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices as literal text (e.g. "10.00", "1.00").
# Assigning a numeric-looking string to .Value on a General-formatted cell
# makes Excel coerce it to a real number (dropping the trailing zeros), so
# for those specific cells we first mark the cell as Text ("@") to preserve
# the exact source formatting, matching the intended inline-string content.

$ws.Range("D2").Value = '69.433.99'
$ws.Range("D3").Value = '3.689.93'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '685.48'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.20'
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +1.12%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000233'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '4.313.42'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.48'
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D15").Value = '3.694.08'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").Value = '69.398.23'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.85'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("E19").Value = '  -1.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '470.08'
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.00'
$ws.Range("E21").Value = '  +2.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.651'
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.79'
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("D24").Value = '3.836.88'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.01'
$ws.Range("E27").Value = '  -4.78%  '
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.72'
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("E30").Value = '  -3.49%  '
$ws.Range("E31").Value = '  -4.63%  '
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.95'
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").Value = '3.663.47'
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("E36").Value = '  -4.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.18'
$ws.Range("E37").Value = '  -3.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.15'
$ws.Range("E38").Value = '  +1.84%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.24'
$ws.Range("E39").Value = '  +2.82%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.943'
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '165.71'
$ws.Range("E44").Value = '  +3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.54'
$ws.Range("E45").Value = '  -1.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000283'
$ws.Range("E46").Value = '  +3.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.13'
$ws.Range("E47").Value = '  +6.83%  '
$ws.Range("E48").Value = '  -4.77%  '
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.98'
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.79'
$ws.Range("E51").Value = '  -2.52%  '
